# Insert a new worksheet "2022-Q1" before the "总计" sheet, so the final
# sheet order is: 2020-Q4, 2021-Q4, 2022-Q1, 总计.
$wb = $excel.ActiveWorkbook
$beforeSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q1"

# Re-fetch sheets fresh by name (avoid stale references captured before Add()).
$qSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# --- Populate the new "2022-Q1" sheet with fund holdings data ---

# Copy the styled header row (基金代码/基金名称/...) from the "2021-Q4" sheet
# so the new header uses the same bold/centered/bordered style.
$templateSheet.Range("B1:H1").Copy($qSheet.Range("B1:H1"))

# Column A (the numeric row index) also uses that same header style in the
# source data, e.g. on the "2021-Q4" sheet; replicate that for rows 2 and 3.
$templateSheet.Range("A2").Copy($qSheet.Range("A2:A3"))

$qSheet.Cells.Item(1, 2).Value = "基金代码"
$qSheet.Cells.Item(1, 3).Value = "基金名称"
$qSheet.Cells.Item(1, 4).Value = "基金规模"
$qSheet.Cells.Item(1, 5).Value = "股票总仓位"
$qSheet.Cells.Item(1, 6).Value = "仓位占比"
$qSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$qSheet.Cells.Item(1, 8).Value = "仓位排名"

# Force the numeric-looking text fields (fund code, size, position figures) to
# be stored as text, matching the source data, by pre-setting the number
# format to Text before assigning the values.
$qSheet.Range("B2:B3").NumberFormat = "@"
$qSheet.Range("D2:G3").NumberFormat = "@"

$qSheet.Cells.Item(2, 1).Value = 0
$qSheet.Cells.Item(2, 2).Value = "519656"
$qSheet.Cells.Item(2, 3).Value = "银河灵活配置混合 - A"
$qSheet.Cells.Item(2, 4).Value = "0.72"
$qSheet.Cells.Item(2, 5).Value = "59.27"
$qSheet.Cells.Item(2, 6).Value = "2.99"
$qSheet.Cells.Item(2, 7).Value = "0.0215"
$qSheet.Cells.Item(2, 8).Value = 8

$qSheet.Cells.Item(3, 1).Value = 1
$qSheet.Cells.Item(3, 2).Value = "519657"
$qSheet.Cells.Item(3, 3).Value = "银河灵活配置混合 - C"
$qSheet.Cells.Item(3, 4).Value = "0.33"
$qSheet.Cells.Item(3, 5).Value = "59.27"
$qSheet.Cells.Item(3, 6).Value = "2.99"
$qSheet.Cells.Item(3, 7).Value = "0.0099"
$qSheet.Cells.Item(3, 8).Value = 8

# --- Update the "总计" (totals) sheet: insert a new first data row for 2022-Q1,
#     pushing the existing 2021-Q4 / 2020-Q4 rows down by one. ---
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0.07000000000000001

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 6
$totalSheet.Cells.Item(3, 4).Value = 0.25

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.03

# Row 4 (2020-Q4) is brand new, so it doesn't inherit the "s=2" style that
# column A carries on the other rows. Copy that format over explicitly.
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(4, 1).PasteSpecial(-4122)
